$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "27.382.94"
Set-TextValue $ws "E2" "  -1.68%  "
Set-TextValue $ws "D3" "1.729.69"
Set-TextValue $ws "E3" "  -1.88%  "
Set-TextValue $ws "D4" "1.004"
Set-TextValue $ws "E4" "  +0.15%  "
Set-TextValue $ws "D5" "322.34"
Set-TextValue $ws "E5" "  -0.22%  "
Set-TextValue $ws "D6" "1.003"
Set-TextValue $ws "E6" "  +0.17%  "
Set-TextValue $ws "D7" "0.4522"
Set-TextValue $ws "E7" "  +5.84%  "
Set-TextValue $ws "D8" "0.3513"
Set-TextValue $ws "E8" "  -3.13%  "
Set-TextValue $ws "B9" "Dogecoin"
Set-TextValue $ws "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws "D9" "0.07331"
Set-TextValue $ws "E9" "  -3.27%  "
Set-TextValue $ws "B10" "OKB"
Set-TextValue $ws "C10" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D10" "41.49"
Set-TextValue $ws "E10" "  -2.94%  "
Set-TextValue $ws "D11" "1.074"
Set-TextValue $ws "E11" "  -2.06%  "
Set-TextValue $ws "D12" "1.003"
Set-TextValue $ws "E12" "  +0.19%  "
Set-TextValue $ws "D13" "20.34"
Set-TextValue $ws "E13" "  -2.06%  "
Set-TextValue $ws "D14" "5.898"
Set-TextValue $ws "E14" "  -2.86%  "
Set-TextValue $ws "D15" "7.031"
Set-TextValue $ws "E15" "  -3.43%  "
Set-TextValue $ws "D16" "1.731.21"
Set-TextValue $ws "E16" "  -1.27%  "
Set-TextValue $ws "D17" "91.09"
Set-TextValue $ws "E17" "  -0.47%  "
Set-TextValue $ws "D18" "0.00001049"
Set-TextValue $ws "E18" "  -1.81%  "
Set-TextValue $ws "E19" "  -0.87%  "
Set-TextValue $ws "D20" "1.003"
Set-TextValue $ws "E20" "  +0.24%  "
Set-TextValue $ws "D21" "16.55"
Set-TextValue $ws "E21" "  -3.15%  "
Set-TextValue $ws "D22" "5.723"
Set-TextValue $ws "E22" "  -3.23%  "
Set-TextValue $ws "D23" "27.437.06"
Set-TextValue $ws "E23" "  -1.55%  "
Set-TextValue $ws "D24" "11.03"
Set-TextValue $ws "E24" "  -2.00%  "
Set-TextValue $ws "D25" "2.076"
Set-TextValue $ws "E25" "  -2.28%  "
Set-TextValue $ws "D26" "161.79"
Set-TextValue $ws "E26" "  +1.01%  "
Set-TextValue $ws "D27" "19.79"
Set-TextValue $ws "E27" "  -2.53%  "
Set-TextValue $ws "D28" "1.924.92"
Set-TextValue $ws "E28" "  -1.61%  "
Set-TextValue $ws "D29" "2.044"
Set-TextValue $ws "E29" "  -4.69%  "
Set-TextValue $ws "D30" "124.32"
Set-TextValue $ws "E30" "  -0.65%  "
Set-TextValue $ws "D31" "1.042"
Set-TextValue $ws "E31" "  -7.19%  "
Set-TextValue $ws "D32" "0.09092"
Set-TextValue $ws "E32" "  +1.81%  "
Set-TextValue $ws "D33" "3.653"
Set-TextValue $ws "E33" "  -0.91%  "
Set-TextValue $ws "D34" "5.330"
Set-TextValue $ws "E34" "  -4.48%  "
Set-TextValue $ws "B35" "VeChain"
Set-TextValue $ws "C35" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D35" "0.02257"
Set-TextValue $ws "E35" "  -2.14%  "
Set-TextValue $ws "B36" "Aptos"
Set-TextValue $ws "C36" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D36" "11.59"
Set-TextValue $ws "E36" "  -5.53%  "
Set-TextValue $ws "D37" "0.05941"
Set-TextValue $ws "E37" "  -1.64%  "
Set-TextValue $ws "D38" "0.2044"
Set-TextValue $ws "E38" "  -3.27%  "
Set-TextValue $ws "D39" "0.6204"
Set-TextValue $ws "E39" "  -2.56%  "
Set-TextValue $ws "D40" "4.846"
Set-TextValue $ws "E40" "  -3.00%  "
Set-TextValue $ws "D41" "1.185"
Set-TextValue $ws "E41" "  +0.23%  "
Set-TextValue $ws "D42" "1.367"
Set-TextValue $ws "E42" "  -2.37%  "
Set-TextValue $ws "D43" "7.678"
Set-TextValue $ws "E43" "  -3.11%  "
Set-TextValue $ws "D44" "13.04"
Set-TextValue $ws "E44" "  -2.38%  "
Set-TextValue $ws "D45" "3.690"
Set-TextValue $ws "E45" "  -0.20%  "
Set-TextValue $ws "D46" "0.5777"
Set-TextValue $ws "E46" "  -1.93%  "
Set-TextValue $ws "D47" "121.65"
Set-TextValue $ws "E47" "  -1.29%  "
Set-TextValue $ws "D48" "1.914"
Set-TextValue $ws "E48" "  -4.15%  "
Set-TextValue $ws "D49" "0.06822"
Set-TextValue $ws "E49" "  -0.25%  "
Set-TextValue $ws "D50" "1.106"
Set-TextValue $ws "E50" "  -6.89%  "
Set-TextValue $ws "D51" "70.76"
Set-TextValue $ws "E51" "  -4.30%  "
